$d = $word.ActiveDocument

$replacements = @(
    @("56×19=1064", "27×26=702"),
    @("14×21=294", "59×39=2301"),
    @("72×38=2736", "65×34=2210"),
    @("36×82=2952", "33×20=660"),
    @("11×44=484", "29×40=1160"),
    @("17×85=1445", "60×83=4980"),
    @("14×44=616", "89×77=6853"),
    @("91×40=3640", "97×44=4268"),
    @("62×32=1984", "35×42=1470"),
    @("35×21=735", "76×70=5320"),
    @("90×23=2070", "78×58=4524"),
    @("75×35=2625", "37×77=2849"),
    @("85×85=7225", "11×42=462"),
    @("20×96=1920", "67×62=4154"),
    @("98×66=6468", "42×58=2436"),
    @("93×98=9114", "11×41=451"),
    @("12×34=408", "86×50=4300"),
    @("48×61=2928", "49×83=4067"),
    @("20×93=1860", "51×70=3570"),
    @("47×23=1081", "17×82=1394"),
    @("43×41=1763", "39×84=3276"),
    @("39×97=3783", "42×63=2646"),
    @("13×80=1040", "53×69=3657"),
    @("62×21=1302", "92×46=4232"),
    @("28×61=1708", "44×81=3564")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
